$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.744.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.38%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.727.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9976"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'240.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.84%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9982"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.14%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4835"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.94%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06187"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.723.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.02%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'15.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.41%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.06872"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.6045"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.20%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.472"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.08%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'77.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.9986"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'26.561.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.09%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007155"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.56%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'11.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.948.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.03%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.421"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.550"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.63%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.061"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.51%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'139.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +3.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'106.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.11%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.367"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.014"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.07938"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.672"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.598"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.6185"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'MXToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'2.456"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.02%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'RenderToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.92%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.9974"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.01499"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.73%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.613"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.68%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.91%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.3834"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.26%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'6.792"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.1156"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05360"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.17%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.900"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.86%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'30.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.243"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.77%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'51.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.89%  "
$ws.Range("E51").Style = "Normal"
